$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append three new timesheet rows (science + nature data)
$ws.Range("A28").Value = 45191
$ws.Range("B28").Value = 1
$ws.Range("A29").Value = 45192
$ws.Range("B29").Value = 1
$ws.Range("A30").Value = 45193
$ws.Range("B30").Value = 1

$ws.Range("A27").Copy()
$ws.Range("A28:A30").PasteSpecial(-4122)

# Move view/selection to match the new bottom of the list
$ws.Range("B30").Select()
$excel.ActiveWindow.ScrollRow = 14
